$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# INTERNAL_DEMAND_FACTOR: 0.04 -> 0.02 (drives the balance recalculation)
$ws.Range("M2").Value = 0.02

# Move the active selection to L2 (FOREIGN_SUPPLY_DEMAND_FACTOR)
$ws.Range("L2").Select()
